# commit by chitra 21-09-2022
# Rework the "HomeSearch" lookup sheet:
#  - shared strings "bag"/"saree"/"kurtis" replaced by "print"/"testcase2"/"testcase3"
#  - those replacement strings now live in column B (rows 1-3) instead of column A (rows 4-6)
#  - rows 4-6 are deleted outright (so A4's old "bag" text and rows 5/6 all disappear),
#    then row 4 is rebuilt with just the empty-but-styled B4 cell
#  - selection moves to A4:A6 (active cell A4)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop rows 4-6 entirely (shrinks the sheet dimension to A1:C4; row 4 gets rebuilt below).
$ws.Range("A4:A6").EntireRow.Delete() | Out-Null

# Row 4 keeps its formatted-but-empty B4 cell, but loses A4 ("bag") altogether.
# Re-create B4's formatting by copying B3's (Hyperlink/text) format onto it instead of
# assigning .Style directly, so we reuse the existing cellXfs entry instead of minting a new one.
$ws.Range("B3").Copy() | Out-Null
$ws.Range("B4").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# New header values for column B.
$ws.Range("B1").Value = "print"
$ws.Range("B2").Value = "testcase2"
$ws.Range("B3").Value = "testcase3"

# Match the updated selection recorded in the sheet view.
$ws.Range("A4:A6").Select() | Out-Null
